$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C12) from 2023-10-08 (45207) to 2023-10-09 (45208)
for ($row = 2; $row -le 12; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value = 45208
    }
}
